# edit.ps1 - apply the Report.docx changes described by the commit diff.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Problem Description" heading: the text was split across two runs
#    ("P" + "roblem Description"). Re-assert the text as a single run by
#    doing a no-op Find/Replace over the full heading text - the engine
#    collapses the matched runs into one while keeping the formatting.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Problem Description", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Problem Description", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Report Audience and Stakeholders" heading gains underline formatting
#    (both the run and the paragraph mark), and the _GoBack bookmark is
#    relocated so it wraps this heading's text instead of sitting at the
#    end of the following paragraph.
# ---------------------------------------------------------------------------
$headingPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Report Audience and Stakeholders*") {
        $headingPara = $p
    }
}
# Underline the run and the paragraph mark (pPr/rPr) together.
$headingPara.Range.Font.Underline = 1

# Move the (hidden) _GoBack bookmark from wherever it currently sits onto
# the heading text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$bookmarkRange = $d.Range($headingPara.Range.Start, $headingPara.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# ---------------------------------------------------------------------------
# 3) Remove the (now empty of bookmark) page-break-only paragraph that used
#    to follow the "Report Audience and Stakeholders" section - it held a
#    lastRenderedPageBreak + manual page break and nothing else.
# ---------------------------------------------------------------------------
$pageBreakPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.Length -le 3 -and $p.Range.Bold) {
        $hasBreak = $false
        foreach ($ch in $t.ToCharArray()) {
            if ([int][char]$ch -eq 12 -or [int][char]$ch -eq 13 -or [int][char]$ch -eq 7) {
                $hasBreak = $true
            }
        }
        if ($hasBreak) {
            $pageBreakPara = $p
        }
    }
}
if ($pageBreakPara -ne $null) {
    $pageBreakPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 4) Foursquare developer hyperlink text was split across three runs
#    ("https://de" + "v" + "eloper.foursquare.com/"). Re-set the display
#    text on the Hyperlink object itself so it collapses into a single run
#    while the Hyperlink rStyle / formatting survives.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks.Item($i)
    if ($h.TextToDisplay -eq "https://developer.foursquare.com/") {
        $h.TextToDisplay = "https://developer.foursquare.com/"
    }
}
